# "Generate Report for Handback"
# Localization status report: mark rows as handed back, stamp handback
# datetimes, and fill in the "Latest Target File" / "Latest Handback File"
# columns (with hyperlinks) for both the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$urlMd1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a413b0db8737cd208cc9b62e374f405987781157/e2e/1bd2ba61-5ad7-4d99-b882-fc16d2019343.md"
$urlMd2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a413b0db8737cd208cc9b62e374f405987781157/e2e/f79a32b9-0e45-419f-b9f3-90a59f4ef055.md"
$name1  = "1bd2ba61-5ad7-4d99-b882-fc16d2019343.md"
$name2  = "f79a32b9-0e45-419f-b9f3-90a59f4ef055.md"

# ---------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    (Overview status columns for both locales, and the Status column on
#    each locale sheet.)
# ---------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# 2. Fill "Latest Target File" (I) and "Latest Handback File" (J) columns,
#    plus stamp "Latest Handback DateTime" (K), for each locale sheet.
# ---------------------------------------------------------------------

# --- zh-cn ---
$wsZhCn.Range("I2").Value = $name1
$wsZhCn.Range("J2").Value = "1bd2ba61-5ad7-4d99-b882-fc16d2019343.588d44d748633d42bab33c756ebe4b0d6bd41738.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-17 12:45:04"

$wsZhCn.Range("I3").Value = $name2
$wsZhCn.Range("J3").Value = "f79a32b9-0e45-419f-b9f3-90a59f4ef055.3061dc6e83076380c9a0ba639082cf9376e9a6da.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-08-17 12:45:04"

# --- de-de ---
$wsDeDe.Range("I2").Value = $name1
$wsDeDe.Range("J2").Value = "1bd2ba61-5ad7-4d99-b882-fc16d2019343.588d44d748633d42bab33c756ebe4b0d6bd41738.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-17 12:45:19"

$wsDeDe.Range("I3").Value = $name2
$wsDeDe.Range("J3").Value = "f79a32b9-0e45-419f-b9f3-90a59f4ef055.3061dc6e83076380c9a0ba639082cf9376e9a6da.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-08-17 12:45:19"

# ---------------------------------------------------------------------
# 3. Re-create the hyperlinks on both sheets so that the "Latest Target
#    File" cells (I2/I3) link to the same source docs as A2/A3, in the
#    order A2, I2, A3, I3 (matches the handoff link + new handback link
#    per row).
# ---------------------------------------------------------------------
foreach ($ws in @($wsZhCn, $wsDeDe)) {
    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), $urlMd1, [Type]::Missing, [Type]::Missing, $name1)
    $ws.Hyperlinks.Add($ws.Range("I2"), $urlMd1, [Type]::Missing, [Type]::Missing, $name1)
    $ws.Hyperlinks.Add($ws.Range("A3"), $urlMd2, [Type]::Missing, [Type]::Missing, $name2)
    $ws.Hyperlinks.Add($ws.Range("I3"), $urlMd2, [Type]::Missing, [Type]::Missing, $name2)

    # Restore the existing "HyperLink" cell style on the linked cells
    # (Hyperlinks.Add otherwise mints a brand-new style).
    $ws.Range("A2").Style = "HyperLink"
    $ws.Range("I2").Style = "HyperLink"
    $ws.Range("A3").Style = "HyperLink"
    $ws.Range("I3").Style = "HyperLink"
}

# ---------------------------------------------------------------------
# 4. Column width adjustments (widen the Status / Target / Handback
#    columns to fit the new, longer text).
# ---------------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 29.9777047293527
$wsOverview.Columns.Item(6).ColumnWidth = 29.9777047293527

foreach ($ws in @($wsZhCn, $wsDeDe)) {
    $ws.Columns.Item(3).ColumnWidth = 29.9777047293527
    $ws.Columns.Item(9).ColumnWidth = 40
    $ws.Columns.Item(10).ColumnWidth = 40
}
